# Update column G ("K") values on Sheet1 per regenerated save_data.
# Commit message: "regen save_data to use K instead of Strike#, regen std/mean,
# calc and write s_vals" -- the K column (formerly a raw strike count) is
# recomputed to a smaller per-game value for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newK = @{
    2  = 2
    3  = 2
    4  = 1
    6  = 2
    7  = 0
    8  = 1
    9  = 2
    10 = 0
    11 = 2
    13 = 1
    14 = 1
    15 = 0
    16 = 0
    17 = 1
    18 = 1
    19 = 2
    20 = 0
    21 = 0
    22 = 1
    23 = 0
    24 = 1
    25 = 2
    26 = 0
    27 = 1
    28 = 1
    29 = 2
    30 = 0
    31 = 0
    32 = 0
    33 = 5
    34 = 1
    35 = 0
    36 = 0
    37 = 2
    39 = 2
    40 = 1
    41 = 1
    43 = 2
    44 = 0
    45 = 1
    46 = 1
    47 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
